# Restructuring Manage view: 50% done
#
# Adds three new localization key/value rows (Day, Days, AddToTop) to the
# bottom of both the "en" (sheet1) and "de" (sheet2) worksheets, each of
# which is a simple two-column key -> translated-value lookup table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "en"
$ws2 = $wb.Worksheets.Item(2)   # "de"

# The last populated row in both sheets is 225; we copy its formatting
# (via PasteSpecial formats-only) onto each freshly used cell before
# writing its value, so the new rows inherit the existing cell styles
# (A column: wrap + vertical-center, B column: number-format/wrap) without
# Excel minting brand-new style records.

# --- Row 226: Day / Tag ---------------------------------------------------
$ws1.Range("A225").Copy()
$ws1.Range("A226").PasteSpecial(-4122)
$ws1.Range("A226").Value() = "Day"
$ws1.Range("B225").Copy()
$ws1.Range("B226").PasteSpecial(-4122)
$ws1.Range("B226").Value() = "Day"

$ws2.Range("A225").Copy()
$ws2.Range("A226").PasteSpecial(-4122)
$ws2.Range("A226").Value() = "Day"
$ws2.Range("B225").Copy()
$ws2.Range("B226").PasteSpecial(-4122)
$ws2.Range("B226").Value() = "Tag"

# --- Row 227: Days / Tage -------------------------------------------------
$ws1.Range("A225").Copy()
$ws1.Range("A227").PasteSpecial(-4122)
$ws1.Range("A227").Value() = "Days"
$ws1.Range("B225").Copy()
$ws1.Range("B227").PasteSpecial(-4122)
$ws1.Range("B227").Value() = "Days"

$ws2.Range("A225").Copy()
$ws2.Range("A227").PasteSpecial(-4122)
$ws2.Range("A227").Value() = "Days"
$ws2.Range("B225").Copy()
$ws2.Range("B227").PasteSpecial(-4122)
$ws2.Range("B227").Value() = "Tage"

# --- Row 228: AddToTop / Add to top / Nach oben erweitern -----------------
$ws1.Range("A225").Copy()
$ws1.Range("A228").PasteSpecial(-4122)
$ws1.Range("A228").Value() = "AddToTop"
$ws1.Range("B225").Copy()
$ws1.Range("B228").PasteSpecial(-4122)
$ws1.Range("B228").Value() = "Add to top"

$ws2.Range("A225").Copy()
$ws2.Range("A228").PasteSpecial(-4122)
$ws2.Range("A228").Value() = "AddToTop"
$ws2.Range("B225").Copy()
$ws2.Range("B228").PasteSpecial(-4122)
$ws2.Range("B228").Value() = "Nach oben erweitern"

# Update the selected cell shown in each sheet's view to reflect the new
# bottom of the table (select sheet2's target last so it keeps being the
# active/tabSelected sheet, matching the source workbook).
[void]$ws1.Range("A228").Select()
[void]$ws2.Range("B229").Select()
